$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "POINT_C_GACHA_NORMAL"
$ws.Range("A7").Value = "POINT_C_GACHA_SPECIAL"
$ws.Range("A8").Value = "POINT_C_GACHA_DESTINY"

$ws.Range("B6").Value = "일반 가챠 재화"
$ws.Range("B7").Value = "슾셜 가챠 재화"
$ws.Range("B8").Value = "운명 가챠 재화"

$ws.Range("C6:C7").Select() | Out-Null
